$d = $word.ActiveDocument

# Locate the last paragraph in the document ("Iterations: 989") and create a
# fresh trailing paragraph mark after it so we have a safe insertion point
# that is not the document's final mark (Word will not let us delete that
# one directly).
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$trailingPara = $d.Paragraphs.Last
$insertionRange = $trailingPara.Range

$xmlPayload = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Comment</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p><w:p><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="新宋体" w:eastAsia="新宋体" w:cs="新宋体"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="新宋体" w:eastAsia="新宋体" w:cs="新宋体"/><w:color w:val="008000"/><w:kern w:val="0"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>// 1. The residual-based and consecutive approximation stopping conditions result in identical solution and quite similar</w:t></w:r></w:p><w:p><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="新宋体" w:eastAsia="新宋体" w:cs="新宋体"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="新宋体" w:eastAsia="新宋体" w:cs="新宋体"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve">    </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="新宋体" w:eastAsia="新宋体" w:cs="新宋体"/><w:color w:val="008000"/><w:kern w:val="0"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>// convergence speed.</w:t></w:r></w:p><w:p><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="新宋体" w:eastAsia="新宋体" w:cs="新宋体"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="新宋体" w:eastAsia="新宋体" w:cs="新宋体"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve">    </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="新宋体" w:eastAsia="新宋体" w:cs="新宋体"/><w:color w:val="008000"/><w:kern w:val="0"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>// 2. Compare the convergence speed of the Jacobi method, the GS method and the SOR method:</w:t></w:r></w:p><w:p><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="新宋体" w:eastAsia="新宋体" w:cs="新宋体"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="新宋体" w:eastAsia="新宋体" w:cs="新宋体"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve">    </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="新宋体" w:eastAsia="新宋体" w:cs="新宋体"/><w:color w:val="008000"/><w:kern w:val="0"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>// SOR with a good choice of omega &gt; GS &gt; Jacobi &gt; SOR with a bad choice of omega</w:t></w:r></w:p><w:p><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:jc w:val="left"/><w:rPr><w:rFonts w:ascii="新宋体" w:eastAsia="新宋体" w:cs="新宋体"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="新宋体" w:eastAsia="新宋体" w:cs="新宋体"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve">    </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="新宋体" w:eastAsia="新宋体" w:cs="新宋体"/><w:color w:val="008000"/><w:kern w:val="0"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>// 3. For the SOR method, when omega is near to 0 or 2, the convergence speed is rather low compared to the omega close to 1</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="新宋体" w:eastAsia="新宋体" w:cs="新宋体"/><w:color w:val="000000"/><w:kern w:val="0"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve">    </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="新宋体" w:eastAsia="新宋体" w:cs="新宋体"/><w:color w:val="008000"/><w:kern w:val="0"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t>// 4. For the SOR method, when omega is near to 1, the solution and the convergence speed is similar to GS method (by definition).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionRange.InsertXML($xmlPayload)

# InsertXML pushed all of our new paragraphs in just before the now-stale
# trailing paragraph mark that we created above, leaving an extra blank
# paragraph at the very end of the document. Remove that stray paragraph
# mark (together with the mark that currently ends our last inserted
# paragraph) so the document text naturally regains a single closing mark
# -- this preserves the pPr/formatting of our last inserted paragraph
# (deleting only the boundary mark, rather than the whole pair, would
# instead make the merged paragraph inherit the stray paragraph's empty
# formatting).
$paras = $d.Paragraphs
$n = $paras.Count
$lastContentPara = $paras.Item($n - 1)
$staleTrailingPara = $paras.Item($n)
$seam = $d.Range($lastContentPara.Range.End - 1, $staleTrailingPara.Range.End)
$seam.Delete()

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
